$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "15-24"
$ws.Range("C2").Value = "15-24"
$ws.Range("D2").Value = "7-16"
$ws.Range("G2").Value = "10-19"
$ws.Range("D3").Value = "15-24"
$ws.Range("D4").Value = "15-24"
$ws.Range("G4").Value = "10-19"
$ws.Range("B5").Value = "15-24"
$ws.Range("D5").Value = "15-24"
$ws.Range("G5").Value = "10-19"
$ws.Range("D6").Value = "15-24"
$ws.Range("D9").Value = "15-24"
$ws.Range("F9").Value = "10-19"
$ws.Range("B10").Value = "15-24"
$ws.Range("C10").Value = "7-16"
$ws.Range("F10").Value = "10-19"
$ws.Range("B11").Value = "15-24"
$ws.Range("C11").Value = "7-16"
$ws.Range("F11").Value = "15-24"
$ws.Range("G11").Value = "10-19"
$ws.Range("C12").Value = "15-24"
$ws.Range("B13").Value = "15-24"
$ws.Range("C13").Value = "15-24"
$ws.Range("G13").Value = "7-16"
$ws.Range("D16").Value = "15-24"
$ws.Range("G16").Value = "7-16"
$ws.Range("C17").Value = "15-24"
$ws.Range("G17").Value = "10-19"
$ws.Range("B18").Value = "15-24"
$ws.Range("G18").Value = "10-19"
$ws.Range("B19").Value = "15-24"
$ws.Range("C19").Value = "15-24"
$ws.Range("D19").Value = "7-16"
$ws.Range("E19").Value = "10-19"
$ws.Range("C20").Value = "15-24"
$ws.Range("E20").Value = "10-19"
$ws.Range("B23").Value = "15-24"
$ws.Range("D23").Value = "10-19"
$ws.Range("F23").Value = "off"
$ws.Range("B24").Value = "15-24"
$ws.Range("F24").Value = "10-19"
$ws.Range("B27").Value = "15-24"
$ws.Range("D27").Value = "off"
$ws.Range("G27").Value = "10-19"
$ws.Range("B30").Value = "15-24"
$ws.Range("D30").Value = "10-19"
$ws.Range("E30").Value = "off"
$ws.Range("B31").Value = "15-24"
$ws.Range("D31").Value = "10-19"
$ws.Range("E31").Value = "off"
